$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("L T", " daihostvlg", " 2024-05-07 18:18:52", " Tài khoản đã mở cước - Không thể thực hiện mở cước", "Group: Test"),
    @("L T", " daihostvlg ", " 2024-05-07 19:47:57", " Tài khoản đã mở cước - Không thể thực hiện mở cước", "Chat trực tiếp với Bot"),
    @("L T", " daihostvlg", " 2024-05-07 20:49:40", " Tài khoản đã mở cước - Không thể thực hiện mở cước", "Group: Test")
)

$startRow = 155
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}
